$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new values (row 6 deleted below) ---
$data = New-Object "object[,]" 4,34
$data[0,0] = 45152.50694444445
$data[0,1] = 22.1
$data[0,2] = 15.162
$data[0,3] = 4.209
$data[0,4] = 46.427
$data[0,5] = 38.398
$data[0,6] = 17.391
$data[0,7] = 57.345
$data[0,8] = 26.759
$data[0,9] = 11.337
$data[0,10] = 17.452
$data[0,11] = 18.439
$data[0,12] = 19.302
$data[0,13] = 5.553
$data[0,14] = 17.294
$data[0,15] = 24.297
$data[0,16] = 14.517
$data[0,17] = 3.762
$data[0,18] = 2.458
$data[0,19] = 255.893
$data[0,20] = 48.146
$data[0,21] = 15.963
$data[0,22] = 31.894
$data[0,23] = 16.657
$data[0,24] = 2.022
$data[0,25] = 28.586
$data[0,26] = 14.1
$data[0,27] = 12.689
$data[0,28] = 14.848
$data[0,29] = 19.074
$data[0,30] = 3.64
$data[0,31] = 50.56
$data[0,32] = 8.856
$data[0,33] = 19.957
$data[1,0] = 45152.51388888889
$data[1,1] = 20.178
$data[1,2] = 14.467
$data[1,3] = 1.944
$data[1,4] = 43.239
$data[1,5] = 35.765
$data[1,6] = 15.88
$data[1,7] = 61.805
$data[1,8] = 24.432
$data[1,9] = 10.655
$data[1,10] = 16.046
$data[1,11] = 17.406
$data[1,12] = 18.213
$data[1,13] = 5.073
$data[1,14] = 15.79
$data[1,15] = 22.348
$data[1,16] = 13.415
$data[1,17] = 1.65
$data[1,18] = 1.176
$data[1,19] = 233.037
$data[1,20] = 44.186
$data[1,21] = 14.575
$data[1,22] = 29.463
$data[1,23] = 15.683
$data[1,24] = 1.891
$data[1,25] = 29.869
$data[1,26] = 12.874
$data[1,27] = 11.559
$data[1,28] = 13.553
$data[1,29] = 18.166
$data[1,30] = 1.294
$data[1,31] = 55.879
$data[1,32] = 8.132
$data[1,33] = 18.223
$data[2,0] = 45152.52083333334
$data[2,1] = 7.206
$data[2,2] = 4.959
$data[2,3] = 1.016
$data[2,4] = 15.251
$data[2,5] = 12.618
$data[2,6] = 5.672
$data[2,7] = 26.785
$data[2,8] = 8.726000000000001
$data[2,9] = 3.736
$data[2,10] = 5.586
$data[2,11] = 6.216
$data[2,12] = 6.394
$data[2,13] = 1.817
$data[2,14] = 5.639
$data[2,15] = 7.941
$data[2,16] = 4.958
$data[2,17] = 1.003
$data[2,18] = 0.542
$data[2,19] = 78.54300000000001
$data[2,20] = 15.975
$data[2,21] = 5.205
$data[2,22] = 10.498
$data[2,23] = 5.68
$data[2,24] = 0.598
$data[2,25] = 12.28
$data[2,26] = 4.598
$data[2,27] = 4.229
$data[2,28] = 4.939
$data[2,29] = 6.441
$data[2,30] = 0.784
$data[2,31] = 24.37
$data[2,32] = 2.827
$data[2,33] = 6.509
$data[3,0] = 45152.52777777778
$data[3,1] = 11.05
$data[3,2] = 7.96
$data[3,3] = 0.93
$data[3,4] = 23.72
$data[3,5] = 19.61
$data[3,6] = 8.699999999999999
$data[3,7] = 33.3
$data[3,8] = 13.38
$data[3,9] = 5.84
$data[3,10] = 8.76
$data[3,11] = 9.619999999999999
$data[3,12] = 10.01
$data[3,13] = 2.78
$data[3,14] = 8.65
$data[3,15] = 12.22
$data[3,16] = 7.4
$data[3,17] = 0.79
$data[3,18] = 0.55
$data[3,19] = 124.28
$data[3,20] = 24.14
$data[3,21] = 7.98
$data[3,22] = 16.07
$data[3,23] = 8.65
$data[3,24] = 1.02
$data[3,25] = 15.98
$data[3,26] = 7.05
$data[3,27] = 6.33
$data[3,28] = 7.43
$data[3,29] = 10.04
$data[3,30] = 0.5600000000000001
$data[3,31] = 29.88
$data[3,32] = 4.44
$data[3,33] = 9.98
$ws.Range("A2:AH5").Value = $data

# --- Delete row 6 (reduces data from 5 rows to 4 data rows) ---
$ws.Rows.Item(6).Delete()

# --- Adjust column widths (character-width model: stored = ColumnWidth + 0.8333333333333333) ---
$offset = 0.8333333333333333
$widthCols8 = @(3,10,11,17,27,28,29)
foreach ($c in $widthCols8) {
    $ws.Columns.Item($c).ColumnWidth = 8 - $offset
}
$ws.Columns.Item(20).ColumnWidth = 9 - $offset

Write-Output "Edit complete"
